$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - update "想去人数" (want-to-go count) column F
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value  = 1452
$wsExhibit.Range("F5").Value  = 12125
$wsExhibit.Range("F6").Value  = 4475
$wsExhibit.Range("F9").Value  = 31
$wsExhibit.Range("F11").Value = 2589
$wsExhibit.Range("F12").Value = 1121
$wsExhibit.Range("F13").Value = 183
$wsExhibit.Range("F14").Value = 62
$wsExhibit.Range("F15").Value = 5249
$wsExhibit.Range("F19").Value = 11428
$wsExhibit.Range("F20").Value = 11491
$wsExhibit.Range("F22").Value = 57
$wsExhibit.Range("F25").Value = 55

# Sheet "全部类型" (all types) - same updates, rows shifted by the extra
# "苏州·龙猫和他的朋友·动漫作品音乐会" entry present in this sheet
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value  = 1452
$wsAll.Range("F5").Value  = 12125
$wsAll.Range("F6").Value  = 4475
$wsAll.Range("F9").Value  = 31
$wsAll.Range("F11").Value = 2589
$wsAll.Range("F13").Value = 1121
$wsAll.Range("F14").Value = 183
$wsAll.Range("F15").Value = 62
$wsAll.Range("F16").Value = 5249
$wsAll.Range("F20").Value = 11428
$wsAll.Range("F21").Value = 11491
$wsAll.Range("F23").Value = 57
$wsAll.Range("F26").Value = 55
